# Add a new sold item row ("CETAL 100MG/ML ORAL DROPS 15 ML") to the day-sale
# report as the 3rd product line (pushing IVYMOND SYRUP, TIRATAM 500MG 30 F.C.
# TABLETS and the syringes row down by one row each), refresh the running
# total, and bump the generated-at timestamp in the footer by one minute.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlShiftDown   = -4121
$xlPasteFormats = -4122

# --- 1. Insert a blank row at row 9, shifting the existing rows 9-13 down
#        to 10-14 (data rows, the totals row and the footer row alike). ---
$ws.Rows("9:9").Insert($xlShiftDown)

# --- 2. Clone the look of an existing "product" row (row 10, which used to
#        be row 9 before the insert) onto the new row 9: number formats,
#        fonts, borders, fill, alignment, etc. ---
$ws.Range("A10:Q10").Copy()
$ws.Range("A9:Q9").PasteSpecial($xlPasteFormats)
$ws.Application.CutCopyMode = $false

# Row height matches the other product rows (same as the row it was cloned from).
$ws.Rows("9:9").RowHeight = 25.5

# Recreate the merged cell layout used by every product row.
$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

# --- 3. Fill in the new product's data. The number formats were already
#        brought over by the formats paste above. H9/N9/Q9 already use a
#        text ("@") format so plain assignment keeps them as text; L9/P9 use
#        a numeric format (the report still stores their value as text, same
#        as every other product row), so flip to text just long enough to
#        assign the value without leaving a literal apostrophe behind, then
#        restore the original numeric format. ---
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "CETAL 100MG/ML ORAL DROPS 15 ML"
$ws.Range("H9").Value = "6:0"
$ws.Range("N9").Value = "23.00"
$ws.Range("Q9").Value = "1:0"

$fmtL = $ws.Range("L9").NumberFormat
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "1"
$ws.Range("L9").NumberFormat = $fmtL

$fmtP = $ws.Range("P9").NumberFormat
$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value = "23.0000"
$ws.Range("P9").NumberFormat = $fmtP

# --- 4. Renumber the product rows that got pushed down. ---
$ws.Range("A10").Value = 4
$ws.Range("A11").Value = 5
$ws.Range("A12").Value = 6

# --- 5. Refresh the running total (now on row 13 after the shift). ---
$ws.Range("P13").Value = 336.2
$ws.Rows("13:13").RowHeight = 24.75

# --- 6. Bump the generated-at timestamp in the footer (now row 14). ---
$ws.Range("A14").Value = "Wednesday, 10 September, 2025 9:42 AM"
